# Intraday standard output (Mar-Sab) update
# Applies the per-row intraday metric updates for 2026-01-16 workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------- Row 2 (AUID) ----------
$ws.Range("K2").Value = 2.02
$ws.Range("L2").Value = 1.92
$ws.Range("P2").Value = 1177041899
$ws.Range("S2").Value = 1.06
$ws.Range("T2").Value = 2.49
$ws.Range("U2").Value = 1.02
$ws.Range("V2").Value = 1.95
$ws.Range("X2").Value = "2026-01-15 09:07:00"
$ws.Range("Z2").Value = 1.91
$ws.Range("AA2").Value = 61318474
$ws.Range("AD2").Value = 66890046
$ws.Range("AG2").Value = 156811048
$ws.Range("AJ2").Value = 343319388
$ws.Range("AN2").Value = 355275705
$ws.Range("AR2").Value = 364182860
$ws.Range("AV2").Value = 510562949

# ---------- Row 3 (BNKK) ----------
$ws.Range("K3").Value = 4.17
$ws.Range("L3").Value = 4.13
$ws.Range("P3").Value = 259786540
$ws.Range("S3").Value = 2.79
$ws.Range("T3").Value = 6.07
$ws.Range("U3").Value = 2.79
$ws.Range("V3").Value = 4.15
$ws.Range("X3").Value = "2026-01-15 08:18:00"
$ws.Range("Y3").Value = 4.27
$ws.Range("AA3").Value = 35061826
$ws.Range("AD3").Value = 36214438
$ws.Range("AG3").Value = 39680986
$ws.Range("AJ3").Value = 41685503
$ws.Range("AN3").Value = 47185474
$ws.Range("AR3").Value = 48469507
$ws.Range("AV3").Value = 187205552

# ---------- Row 4 (CGTL) ----------
$ws.Range("L4").Value = 3.4
$ws.Range("T4").Value = 4.25
$ws.Range("U4").Value = 2.41
$ws.Range("V4").Value = 3.36
$ws.Range("X4").Value = "2026-01-15 08:38:00"
$ws.Range("Z4").Value = 3.36

# ---------- Row 5 (RILY) ----------
$ws.Range("K5").Value = 9.699999999999999
$ws.Range("L5").Value = 10.04
$ws.Range("T5").Value = 11.65
$ws.Range("U5").Value = 9.43
$ws.Range("V5").Value = 10.01
$ws.Range("X5").Value = "2026-01-15 08:00:00"
$ws.Range("Y5").Value = 10.04
$ws.Range("Z5").Value = 9.380000000000001
$ws.Range("AC5").Value = 9.380000000000001
$ws.Range("AF5").Value = 9.380000000000001
$ws.Range("AI5").Value = 9.380000000000001
$ws.Range("AM5").Value = 9.380000000000001
$ws.Range("AQ5").Value = 9.380000000000001
$ws.Range("AU5").Value = 9.380000000000001

# ---------- Row 6 (SPHL) ----------
$ws.Range("K6").Value = 7.72
$ws.Range("L6").Value = 7.88
$ws.Range("P6").Value = 102336324
$ws.Range("S6").Value = 2.28
$ws.Range("T6").Value = 12.42
$ws.Range("U6").Value = 2.28
$ws.Range("V6").Value = 7.8
$ws.Range("X6").Value = "2026-01-15 07:54:00"
$ws.Range("AA6").Value = 36088844
$ws.Range("AD6").Value = 36705996
$ws.Range("AG6").Value = 67455332
$ws.Range("AJ6").Value = 73348472
$ws.Range("AN6").Value = 78635312
$ws.Range("AR6").Value = 85734532
$ws.Range("AV6").Value = 93950255

Write-Output "Applied intraday updates for rows 2-6"
